# Generate Report for Handoff
# Updates the row for 44e77bcf-982a-426c-a4c7-fdabc76da2d2 across all sheets:
#  - Status changes from "Handed back: in sync with en-US" to "Ready for handoff"
#  - Handoff date/datetime fields are refreshed to the new handoff timestamp

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-25 09:21:19"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-25 09:21:14"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-25 09:21:19"
